$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add prebuilt construction scenario": revise the enhanced-wastage (column C)
# figures for the affected building materials.
$ws.Range("C3").Value = 0.015
$ws.Range("C4").Value = 0.015
$ws.Range("C5").Value = 0.025
$ws.Range("C6").Value = 0.05
$ws.Range("C9").Value = 0.05

# Turn on AutoFilter for the data header row.
$headerRange = $ws.Range("A1:C1")
$headerRange.AutoFilter() | Out-Null

# Register the (hidden) _FilterDatabase defined name that Excel creates
# alongside an AutoFilter, scoped to this worksheet.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:C1"))
$filterName.Visible = $false

# Move the active selection.
$ws.Range("C37").Select() | Out-Null
